$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the values 1-10 under the existing data (rows 12-21) in column A.
$values = 1..10
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Move the selection to A12, matching the post-edit workbook state.
$ws.Range("A12").Select()
